# Apply updated cryptocurrency price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '35.475.14'
$ws.Range('E2').Value = '  -2.70%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '1.971.73'
$ws.Range('E3').Value = '  -4.02%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  +0.09%  '

# Row 5: BNB
$ws.Range('D5').Value = '''243.84'
$ws.Range('E5').Value = '  +0.94%  '

# Row 7: Solana
$ws.Range('D7').Value = '''56.90'
$ws.Range('E7').Value = '  +4.92%  '

# Row 8: USDC
$ws.Range('E8').Value = '  +0.02%  '

# Row 9: OKB
$ws.Range('D9').Value = '''58.79'
$ws.Range('E9').Value = '  +1.06%  '

# Row 10: Cardano
$ws.Range('D10').Value = '''0.357'
$ws.Range('E10').Value = '  +0.42%  '

# Row 11: Dogecoin
$ws.Range('D11').Value = '''0.0729'
$ws.Range('E11').Value = '  -2.42%  '

# Row 12: TRON
$ws.Range('E12').Value = '  -3.30%  '

# Row 13: Polygon
$ws.Range('D13').Value = '''0.939'
$ws.Range('E13').Value = '  +5.13%  '

# Row 14: Chainlink
$ws.Range('D14').Value = '''14.11'
$ws.Range('E14').Value = '  -3.25%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range('D15').Value = '2.260.91'
$ws.Range('E15').Value = '  -3.98%  '

# Row 16: Polkadot
$ws.Range('D16').Value = '''5.22'
$ws.Range('E16').Value = '  -2.09%  '

# Row 17: WrappedEther
$ws.Range('D17').Value = '1.969.68'
$ws.Range('E17').Value = '  -4.04%  '

# Row 18: Avalanche
$ws.Range('D18').Value = '''17.52'
$ws.Range('E18').Value = '  +5.55%  '

# Row 19: WrappedBTC
$ws.Range('D19').Value = '35.387.73'
$ws.Range('E19').Value = '  -2.78%  '

# Row 20: Litecoin
$ws.Range('D20').Value = '''71.18'
$ws.Range('E20').Value = '  -1.01%  '

# Row 21: ShibaInu
$ws.Range('D21').Value = '0.0₃0838'
$ws.Range('E21').Value = '  -1.94%  '

# Row 22: BitcoinCash
$ws.Range('D22').Value = '''232.16'
$ws.Range('E22').Value = '  -2.20%  '

# Row 23: Uniswap
$ws.Range('D23').Value = '''5.11'
$ws.Range('E23').Value = '  -2.32%  '

# Row 24: Dai
$ws.Range('E24').Value = '  -0.08%  '

# Row 25: PancakeSwap
$ws.Range('E25').Value = '  +19.58%  '

# Row 26: Toncoin
$ws.Range('E26').Value = '  -1.82%  '

# Row 27: Monero
$ws.Range('D27').Value = '''163.08'
$ws.Range('E27').Value = '  +0.35%  '

# Row 28: Cosmos
$ws.Range('D28').Value = '''9.04'
$ws.Range('E28').Value = '  -2.84%  '

# Row 29: EthereumClassic
$ws.Range('D29').Value = '''19.19'
$ws.Range('E29').Value = '  -4.30%  '

# Row 30: Stellar
$ws.Range('E30').Value = '  -2.53%  '

# Row 31: Filecoin
$ws.Range('E31').Value = '  -3.96%  '

# Row 32: ImmutableX
$ws.Range('D32').Value = '''1.12'
$ws.Range('E32').Value = '  -3.25%  '

# Row 33: Hedera
$ws.Range('D33').Value = '''0.0589'
$ws.Range('E33').Value = '  -0.29%  '

# Row 34: Kaspa
$ws.Range('D34').Value = '''0.0914'
$ws.Range('E34').Value = '  +10.46%  '

# Row 35: InternetComputer(DFINITY)
$ws.Range('D35').Value = '''4.25'
$ws.Range('E35').Value = '  -4.89%  '

# Row 36: LidoDAOToken
$ws.Range('D36').Value = '''2.34'
$ws.Range('E36').Value = '  +7.77%  '

# Row 37: BinanceUSD
$ws.Range('E37').Value = '  +0.09%  '

# Row 38: WEMIXToken
$ws.Range('E38').Value = '  -5.32%  '

# Row 39: THORChain
$ws.Range('D39').Value = '''5.10'
$ws.Range('E39').Value = '  +5.76%  '

# Row 40: TrustWalletToken
$ws.Range('E40').Value = '  -2.96%  '

# Row 41: HuobiToken
$ws.Range('E41').Value = '  +1.18%  '

# Row 42: VeChain
$ws.Range('E42').Value = '  -1.30%  '

# Row 43: ARBITRUM
$ws.Range('E43').Value = '  -2.42%  '

# Row 44: Aave
$ws.Range('D44').Value = '''91.00'
$ws.Range('E44').Value = '  -2.44%  '

# Row 45: Maker
$ws.Range('D45').Value = '1.375.88'
$ws.Range('E45').Value = '  +0.09%  '

# Row 46: InjectiveProtocol
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').Value = '''0.0882'
$ws.Range('E46').Value = '  -1.57%  '

# Row 47: Cronos
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '''15.84'
$ws.Range('E47').Value = '  +1.30%  '

# Row 48: FraxShare
$ws.Range('D48').Value = '''7.47'
$ws.Range('E48').Value = '  +2.91%  '

# Row 49: MXToken
$ws.Range('E49').Value = '  +1.00%  '

# Row 50: RenderToken
$ws.Range('D50').Value = '''2.24'
$ws.Range('E50').Value = '  -0.46%  '

# Row 51: MultiversX
$ws.Range('D51').Value = '''45.41'
$ws.Range('E51').Value = '  +2.46%  '
